$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update columns D (Fecha), J (Volumen), K (Precio mínimo), L (Precio máximo),
# M (Precio promedio ponderado), P (Precio $/Kg) for rows 2-20.
# This reflects a re-sequencing of weekly records.

$ws.Cells.Item(2, 4).Value = 44832
$ws.Cells.Item(2, 10).Value = 60
$ws.Cells.Item(2, 11).Value = 17000
$ws.Cells.Item(2, 12).Value = 18000
$ws.Cells.Item(2, 13).Value = 17500
$ws.Cells.Item(2, 16).Value = 1346

$ws.Cells.Item(3, 4).Value = 44874
$ws.Cells.Item(3, 10).Value = 30
$ws.Cells.Item(3, 11).Value = 17000
$ws.Cells.Item(3, 12).Value = 17000
$ws.Cells.Item(3, 13).Value = 17000
$ws.Cells.Item(3, 16).Value = 1308

$ws.Cells.Item(4, 4).Value = 44959
$ws.Cells.Item(4, 10).Value = 30
$ws.Cells.Item(4, 11).Value = 19000
$ws.Cells.Item(4, 12).Value = 19000
$ws.Cells.Item(4, 13).Value = 19000
$ws.Cells.Item(4, 16).Value = 1462

$ws.Cells.Item(5, 4).Value = 44839
$ws.Cells.Item(5, 10).Value = 40
$ws.Cells.Item(5, 11).Value = 15000
$ws.Cells.Item(5, 12).Value = 16000
$ws.Cells.Item(5, 13).Value = 15500
$ws.Cells.Item(5, 16).Value = 1192

$ws.Cells.Item(6, 4).Value = 45155
$ws.Cells.Item(6, 10).Value = 30
$ws.Cells.Item(6, 11).Value = 20000
$ws.Cells.Item(6, 12).Value = 20000
$ws.Cells.Item(6, 13).Value = 20000
$ws.Cells.Item(6, 16).Value = 1538

$ws.Cells.Item(7, 4).Value = 44930
$ws.Cells.Item(7, 10).Value = 30
$ws.Cells.Item(7, 11).Value = 17000
$ws.Cells.Item(7, 12).Value = 17000
$ws.Cells.Item(7, 13).Value = 17000
$ws.Cells.Item(7, 16).Value = 1308

$ws.Cells.Item(8, 4).Value = 44894
$ws.Cells.Item(8, 10).Value = 30
$ws.Cells.Item(8, 11).Value = 18000
$ws.Cells.Item(8, 12).Value = 18000
$ws.Cells.Item(8, 13).Value = 18000
$ws.Cells.Item(8, 16).Value = 1385

$ws.Cells.Item(9, 4).Value = 44859
$ws.Cells.Item(9, 10).Value = 30
$ws.Cells.Item(9, 11).Value = 13000
$ws.Cells.Item(9, 12).Value = 13000
$ws.Cells.Item(9, 13).Value = 13000
$ws.Cells.Item(9, 16).Value = 1000

$ws.Cells.Item(10, 4).Value = 44915
$ws.Cells.Item(10, 10).Value = 50
$ws.Cells.Item(10, 11).Value = 18000
$ws.Cells.Item(10, 12).Value = 18000
$ws.Cells.Item(10, 13).Value = 18000
$ws.Cells.Item(10, 16).Value = 1385

$ws.Cells.Item(11, 4).Value = 44880
$ws.Cells.Item(11, 10).Value = 30
$ws.Cells.Item(11, 11).Value = 17000
$ws.Cells.Item(11, 12).Value = 17000
$ws.Cells.Item(11, 13).Value = 17000
$ws.Cells.Item(11, 16).Value = 1308

$ws.Cells.Item(12, 4).Value = 44810
$ws.Cells.Item(12, 10).Value = 40
$ws.Cells.Item(12, 11).Value = 12000
$ws.Cells.Item(12, 12).Value = 13000
$ws.Cells.Item(12, 13).Value = 12500
$ws.Cells.Item(12, 16).Value = 962

$ws.Cells.Item(13, 4).Value = 44868
$ws.Cells.Item(13, 10).Value = 30
$ws.Cells.Item(13, 11).Value = 18000
$ws.Cells.Item(13, 12).Value = 18000
$ws.Cells.Item(13, 13).Value = 18000
$ws.Cells.Item(13, 16).Value = 1385

$ws.Cells.Item(14, 4).Value = 44797
$ws.Cells.Item(14, 10).Value = 60
$ws.Cells.Item(14, 11).Value = 12000
$ws.Cells.Item(14, 12).Value = 13000
$ws.Cells.Item(14, 13).Value = 12500
$ws.Cells.Item(14, 16).Value = 962

$ws.Cells.Item(15, 4).Value = 44943
$ws.Cells.Item(15, 10).Value = 30
$ws.Cells.Item(15, 11).Value = 17000
$ws.Cells.Item(15, 12).Value = 17000
$ws.Cells.Item(15, 13).Value = 17000
$ws.Cells.Item(15, 16).Value = 1308

$ws.Cells.Item(16, 4).Value = 44804
$ws.Cells.Item(16, 10).Value = 40
$ws.Cells.Item(16, 11).Value = 12000
$ws.Cells.Item(16, 12).Value = 13000
$ws.Cells.Item(16, 13).Value = 12500
$ws.Cells.Item(16, 16).Value = 962

$ws.Cells.Item(17, 4).Value = 44841
$ws.Cells.Item(17, 10).Value = 30
$ws.Cells.Item(17, 11).Value = 18000
$ws.Cells.Item(17, 12).Value = 18000
$ws.Cells.Item(17, 13).Value = 18000
$ws.Cells.Item(17, 16).Value = 1385

$ws.Cells.Item(18, 4).Value = 44846
$ws.Cells.Item(18, 10).Value = 30
$ws.Cells.Item(18, 11).Value = 18000
$ws.Cells.Item(18, 12).Value = 18000
$ws.Cells.Item(18, 13).Value = 18000
$ws.Cells.Item(18, 16).Value = 1385

$ws.Cells.Item(19, 4).Value = 44922
$ws.Cells.Item(19, 10).Value = 30
$ws.Cells.Item(19, 11).Value = 17000
$ws.Cells.Item(19, 12).Value = 17000
$ws.Cells.Item(19, 13).Value = 17000
$ws.Cells.Item(19, 16).Value = 1308

$ws.Cells.Item(20, 4).Value = 44895
$ws.Cells.Item(20, 10).Value = 30
$ws.Cells.Item(20, 11).Value = 18000
$ws.Cells.Item(20, 12).Value = 18000
$ws.Cells.Item(20, 13).Value = 18000
$ws.Cells.Item(20, 16).Value = 1385
